$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.348.42"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.529.66"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D9").Value = "2.527.10"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "2.976.88"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "59.277.47"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "2.529.35"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.25%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  +0.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("E33").Value = "  +8.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.827"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "131.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.61%  "
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("E51").Value = "  -0.53%  "
